# Updates cryptos list values (price & 1h volume change) to match the
# latest scrape, and re-sorts/renames a few rows (28-51) whose underlying
# coin ranking shifted by one position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.461.21"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "3.260.15"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'574.94"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'178.87"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.582"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").Value = "3.262.80"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'0.175"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'0.571"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "'45.39"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "'0.0000268"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "'674.69"
$ws.Range("E14").Value = "  +10.34%  "
$ws.Range("D15").Value = "3.806.61"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'8.31"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "67.698.88"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "3.282.31"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'17.31"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'10.73"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "'0.886"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'17.03"
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("D24").Value = "'5.10"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").Value = "'97.89"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "'2.72"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'9.32"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'32.57"
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.36"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.66"
$ws.Range("E31").Value = "  +4.15%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'582.35"
$ws.Range("E32").Value = "  +6.14%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "3.863.60"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'10.80"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.103"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.36"
$ws.Range("E37").Value = "  -9.14%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'55.22"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.129"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.22"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'32.05"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0668"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.329"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0411"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.127"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'1.38"
$ws.Range("E49").Value = "  +9.13%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.49"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'129.79"
$ws.Range("E51").Value = "  +0.74%  "

Write-Host "Updated cryptos list"
